$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'Datos actualizados a 31 de Marzo de 2020 a las 01:50'
$ws.Range("B4").Value = 163195
$ws.Range("C4").Value = 19704
$ws.Range("D4").Value = 5501
$ws.Range("E4").Value = 154546
$ws.Range("G4").Value = 565
$ws.Range("H4").Value = 3148
$ws.Range("B18").Value = 7448
$ws.Range("C18").Value = 1128
$ws.Range("E18").Value = 6266
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 89
$ws.Range("B21").Value = 4630
$ws.Range("C21").Value = 374
$ws.Range("E21").Value = 4347
$ws.Range("G21").Value = 27
$ws.Range("H21").Value = 163
$ws.Range("E44").Value = 1114
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 46
$ws.Range("A46").Value = 'Panama'
$ws.Range("B46").Value = 1075
$ws.Range("C46").Value = 86
$ws.Range("D46").Value = 4
$ws.Range("E46").Value = 1047
$ws.Range("F46").Value = 36
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 24
$ws.Range("A47").Value = 'Mexico'
$ws.Range("B47").Value = 993
$ws.Range("C47").Value = 145
$ws.Range("D47").Value = 35
$ws.Range("E47").Value = 938
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 20
$ws.Range("E52").Value = 769
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 14
$ws.Range("A59").Value = 'Hong Kong'
$ws.Range("B59").Value = 683
$ws.Range("C59").Value = 41
$ws.Range("D59").Value = 118
$ws.Range("E59").Value = 561
$ws.Range("F59").Value = 5
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 4
$ws.Range("A60").Value = 'Egipto'
$ws.Range("B60").Value = 656
$ws.Range("C60").Value = 47
$ws.Range("D60").Value = 150
$ws.Range("E60").Value = 465
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 41
$ws.Range("A75").Value = 'Tunez'
$ws.Range("B75").Value = 362
$ws.Range("C75").Value = 50
$ws.Range("D75").Value = 3
$ws.Range("E75").Value = 349
$ws.Range("F75").Value = 10
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 10
$ws.Range("A76").Value = 'Bulgaria'
$ws.Range("B76").Value = 359
$ws.Range("C76").Value = 13
$ws.Range("D76").Value = 17
$ws.Range("E76").Value = 334
$ws.Range("F76").Value = 13
$ws.Range("H76").Value = 8
$ws.Range("A77").Value = 'Eslovaquia'
$ws.Range("B77").Value = 336
$ws.Range("C77").Value = 22
$ws.Range("D77").Value = 7
$ws.Range("E77").Value = 329
$ws.Range("F77").Value = 1
$ws.Range("H77").Value = 0
$ws.Range("A78").Value = 'Costa Rica'
$ws.Range("B78").Value = 330
$ws.Range("C78").Value = 16
$ws.Range("D78").Value = 4
$ws.Range("E78").Value = 324
$ws.Range("F78").Value = 7
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 2
$ws.Range("B79").Value = 320
$ws.Range("C79").Value = 16
$ws.Range("E79").Value = 319
$ws.Range("A158").Value = 'Dominica'
$ws.Range("C158").Value = 1
$ws.Range("A159").Value = 'Eritrea'
$ws.Range("A160").Value = 'Guinea Ecuatorial'
$ws.Range("E160").Value = 12
$ws.Range("H160").Value = 0
$ws.Range("A161").Value = 'Islas Caimanes'
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 11
$ws.Range("H161").Value = 1
$ws.Range("A162").Value = 'Mongolia'
$ws.Range("B162").Value = 12
$ws.Range("D162").Value = 2
$ws.Range("E162").Value = 10
